# Weekly update: insert two new price records (week of 2023-05-29, serial 45075)
# at the top of the historical table (rows 616-617), shifting the existing
# rows 616-647 down to 618-649.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 616 (shifts 616:647 -> 618:649)
$ws.Rows("616:617").Insert()

# New row 616: Apio, Americana (o), Primera
$ws.Range("A616").Value = 8
$ws.Range("B616").Value = "Terminal La Palmera de La Serena"
$ws.Range("C616").Value = "Coquimbo"
$ws.Range("D616").Value2 = 45075
$ws.Range("E616").Value = 4
$ws.Range("F616").Value = 100112017
$ws.Range("G616").Value = "Apio"
$ws.Range("H616").Value = "Americana (o)"
$ws.Range("I616").Value = "Primera"
$ws.Range("J616").Value = 800
$ws.Range("K616").Value = 8000
$ws.Range("L616").Value = 9000
$ws.Range("M616").Value = 8500
$ws.Range("N616").Value = "$/docena de matas"
$ws.Range("O616").Value = "Provincia del Elquí"
$ws.Range("P616").Value = 1417
$ws.Range("Q616").Value = 6
$ws.Range("R616").Value = "Hortaliza"

# New row 617: Apio, Americana (o), Segunda
$ws.Range("A617").Value = 8
$ws.Range("B617").Value = "Terminal La Palmera de La Serena"
$ws.Range("C617").Value = "Coquimbo"
$ws.Range("D617").Value2 = 45075
$ws.Range("E617").Value = 4
$ws.Range("F617").Value = 100112017
$ws.Range("G617").Value = "Apio"
$ws.Range("H617").Value = "Americana (o)"
$ws.Range("I617").Value = "Segunda"
$ws.Range("J617").Value = 400
$ws.Range("K617").Value = 6000
$ws.Range("L617").Value = 7000
$ws.Range("M617").Value = 6500
$ws.Range("N617").Value = "$/docena de matas"
$ws.Range("O617").Value = "Provincia del Elquí"
$ws.Range("P617").Value = 1083
$ws.Range("Q617").Value = 6
$ws.Range("R617").Value = "Hortaliza"
